$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (existing data rows end at row 11)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Append 3 more rows identical to the existing data rows (mollie_ABC123 | 1445758 | Ben Gortemaker | 1)
for ($i = 1; $i -le 3; $i++) {
    $newRow = $lastRow + $i
    $ws.Cells.Item($newRow, 1).Value = "mollie_ABC123"
    $ws.Cells.Item($newRow, 2).Value = 1445758
    $ws.Cells.Item($newRow, 3).Value = "Ben Gortemaker"
    $ws.Cells.Item($newRow, 4).Value = 1
}
